$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (s=2 currently) - set font color white
$ws.Range("A2:K2").Font.Color = 16777215

# Title (s=1) - remove dedicated 14pt size, set font color white too, to try to force font merge
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2:G2").Font.Color = 16777215
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215

$ws.Range("H3").Value = -47
$ws.Range("I3").Value = "16-Sep-2025"
